$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.412.75'
$ws.Range('E2').Value = '  +1.88%  '
$ws.Range('D3').Value = '2.354.59'
$ws.Range('E3').Value = '  +1.54%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '''521.62'
$ws.Range('E5').Value = '  +0.37%  '
$ws.Range('D6').Value = '''136.39'
$ws.Range('E6').Value = '  +1.68%  '
$ws.Range('E7').Value = '  +0.33%  '
$ws.Range('E8').Value = '  +0.36%  '
$ws.Range('D9').Value = '2.363.48'
$ws.Range('E9').Value = '  +1.03%  '
$ws.Range('E10').Value = '  -0.54%  '
$ws.Range('D11').Value = '''5.45'
$ws.Range('E11').Value = '  +5.15%  '
$ws.Range('E12').Value = '  -1.22%  '
$ws.Range('E13').Value = '  +0.06%  '
$ws.Range('D14').Value = '''24.34'
$ws.Range('E14').Value = '  +0.80%  '
$ws.Range('D15').Value = '2.773.27'
$ws.Range('E15').Value = '  +1.52%  '
$ws.Range('D16').Value = '57.393.23'
$ws.Range('E16').Value = '  +1.68%  '
$ws.Range('E17').Value = '  +0.13%  '
$ws.Range('D18').Value = '2.358.02'
$ws.Range('E18').Value = '  +0.26%  '
$ws.Range('D19').Value = '''10.63'
$ws.Range('E19').Value = '  +0.61%  '
$ws.Range('D20').Value = '''329.61'
$ws.Range('E20').Value = '  +2.54%  '
$ws.Range('D21').Value = '''4.25'
$ws.Range('E21').Value = '  -1.06%  '
$ws.Range('D22').Value = '''6.73'
$ws.Range('E22').Value = '  +1.29%  '
$ws.Range('D23').Value = '''0.998'
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('D24').Value = '''61.47'
$ws.Range('E24').Value = '  +1.07%  '
$ws.Range('D25').Value = '''0.166'
$ws.Range('E25').Value = '  +3.78%  '
$ws.Range('D26').Value = '''1.00'
$ws.Range('E26').Value = '  +0.75%  '
$ws.Range('D27').Value = '''8.31'
$ws.Range('E27').Value = '  +8.21%  '
$ws.Range('E28').Value = '  +9.07%  '
$ws.Range('E29').Value = '  -0.87%  '
$ws.Range('D30').Value = '0.0₃0746'
$ws.Range('E30').Value = '  +1.23%  '
$ws.Range('D31').Value = '''1.71'
$ws.Range('E31').Value = '  +0.28%  '
$ws.Range('D32').Value = '''6.29'
$ws.Range('E32').Value = '  +0.02%  '
$ws.Range('D33').Value = '''18.59'
$ws.Range('E33').Value = '  +0.82%  '
$ws.Range('E34').Value = '  +0.04%  '
$ws.Range('E35').Value = '  +2.86%  '
$ws.Range('D36').Value = '''0.992'
$ws.Range('E36').Value = '  +0.25%  '
$ws.Range('D37').Value = '''0.929'
$ws.Range('E37').Value = '  -0.49%  '
$ws.Range('D38').Value = '''4.05'
$ws.Range('E38').Value = '  +0.47%  '
$ws.Range('E39').Value = '  +3.51%  '
$ws.Range('D40').Value = '''38.57'
$ws.Range('E40').Value = '  +2.79%  '
$ws.Range('D41').Value = '''151.23'
$ws.Range('E41').Value = '  +7.09%  '
$ws.Range('E42').Value = '  +0.04%  '
$ws.Range('D43').Value = '''3.66'
$ws.Range('E43').Value = '  +1.69%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').Value = '''5.30'
$ws.Range('E44').Value = '  +2.39%  '
$ws.Range('B45').Value = 'Bittensor'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D45').Value = '''283.84'
$ws.Range('E45').Value = '  +1.34%  '
$ws.Range('D46').Value = '''0.0939'
$ws.Range('E46').Value = '  +1.00%  '
$ws.Range('D47').Value = '''0.0509'
$ws.Range('E48').Value = '  +0.93%  '
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').Value = '''18.36'
$ws.Range('E49').Value = '  +6.05%  '
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').Value = '''0.0221'
$ws.Range('E50').Value = '  +2.18%  '
$ws.Range('D51').Value = '''17.62'
$ws.Range('E51').Value = '  +3.79%  '
